$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing tracker rows ---

# Row 9 (Dorman - Rotational Product Engineer): got a response, now planning interview
$ws.Range("F9").Value = "Got response"
$ws.Range("G9").Value = "Planning"

# Row 12 (Boston Scientific - Manufacturing Engineer II): declined
$ws.Range("H12").Value = "Declined"

# Row 19 (Boston Scientific - R&D Engineer I): interview cell cleared, declined
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = "Declined"

# --- Add new row 20: Polaris - Test Cell Technician (cover letter sent) ---
$ws.Range("A20").Value = "Polaris"
$ws.Range("B20").Value = "Test Cell Technician"
$ws.Range("C20").HorizontalAlignment = 5
$ws.Range("D20").Value = "Done"
$ws.Range("D20").HorizontalAlignment = -4108
$ws.Range("E20").Value = "Done"
$ws.Range("E20").HorizontalAlignment = -4108
$ws.Range("F20").Value = "Emailed"
$ws.Range("F20").HorizontalAlignment = -4108
$ws.Range("G20").HorizontalAlignment = -4108
$ws.Range("H20").Value = "Declined"
$ws.Range("H20").HorizontalAlignment = -4108

# --- Resize the table / autofilter to include the new row ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:I20"))

# --- Resize conditional formatting ranges to include row 20 ---
function Resize-CF($probeRange, $exactFormula, $newRange) {
    $fc = $ws.Range($probeRange).FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $item = $fc.Item($i)
        if ($item.Formula1 -eq $exactFormula) {
            $item.ModifyAppliesToRange($ws.Range($newRange))
            return
        }
    }
}

Resize-CF "D2:E19" '=NOT(ISERROR(SEARCH("In progress",D2)))' "D2:E20"
Resize-CF "D2:F19" '=LEN(TRIM(D2))=0' "D2:F20"
Resize-CF "D2:G19" '=NOT(ISERROR(SEARCH("Done",D2)))' "D2:G20"
Resize-CF "F2:F19" '=NOT(ISERROR(SEARCH("Called",F2)))' "F2:F20"
Resize-CF "F2:F19" '=NOT(ISERROR(SEARCH("Emailed",F2)))' "F2:F20"
Resize-CF "F2:F19" '=NOT(ISERROR(SEARCH("Got response",F2)))' "F2:F20"
Resize-CF "F2:F19" '=NOT(ISERROR(SEARCH("LInkedIn",F2)))' "F2:F20"
Resize-CF "F2:F19" '=NOT(ISERROR(SEARCH("Need contact",F2)))' "F2:F20"
Resize-CF "G2:G19" '=NOT(ISERROR(SEARCH("Planning",G2)))' "G2:G20"
Resize-CF "H2:H19" '=NOT(ISERROR(SEARCH("Declined",H2)))' "H2:H20"
Resize-CF "H2:H19" '=NOT(ISERROR(SEARCH("Waiting",H2)))' "H2:H20"
Resize-CF "H2:H19" '=NOT(ISERROR(SEARCH("Recieved",H2)))' "H2:H20"

# --- Resize data validation ranges to include row 20 ---
$ws.Range("D2:E20").Validation.Delete()
$ws.Range("D2:E20").Validation.Add(3, 1, 1, '"In progress,Done"')

$ws.Range("H2:H20").Validation.Delete()
$ws.Range("H2:H20").Validation.Add(3, 1, 1, '"Waiting, Received, Declined"')

$ws.Range("F2:F20").Validation.Delete()
$ws.Range("F2:F20").Validation.Add(3, 1, 1, '"Got response, LinkedIn, Emailed, Called"')

$ws.Range("G2:G20").Validation.Delete()
$ws.Range("G2:G20").Validation.Add(3, 1, 1, '"Planning, Done"')

# --- Update selection to where the user left off ---
$ws.Range("A21").Select()

$wb.Save()
